# Updated with round 12
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 45802
    3  = 45801
    4  = 45801
    5  = 45801
    6  = 45800
    7  = 45801
    8  = 45801
    9  = 45799
    10 = 45802
    11 = 45801
    12 = 45802
    13 = 45801
    14 = 45801
    15 = 45800
    16 = 45802
    17 = 45802
    18 = 45802
    19 = 45799
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
